$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.084.39"
$ws.Range("E2").Value = "  +7.15%  "
$ws.Range("D3").Value = "3.018.11"
$ws.Range("E3").Value = "  +4.23%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.87"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.39"
$ws.Range("E6").Value = "  +9.12%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.015.23"
$ws.Range("E8").Value = "  +4.19%  "
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.03"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +6.93%  "
$ws.Range("E12").Value = "  +5.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  +9.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.48"
$ws.Range("E14").Value = "  +8.46%  "
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "66.015.25"
$ws.Range("D17").Value = "3.518.65"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.95"
$ws.Range("E18").Value = "  +6.27%  "
$ws.Range("D19").Value = "3.011.90"
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.20"
$ws.Range("E20").Value = "  +7.86%  "
$ws.Range("E21").Value = "  +6.27%  "
$ws.Range("E22").Value = "  +4.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.38"
$ws.Range("E23").Value = "  +8.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.21"
$ws.Range("E24").Value = "  +3.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.52"
$ws.Range("E25").Value = "  +5.40%  "
$ws.Range("E26").Value = "  +11.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.76"
$ws.Range("E27").Value = "  +8.20%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.98"
$ws.Range("E29").Value = "  +13.69%  "
$ws.Range("E30").Value = "  +17.75%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +5.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").Value = "  +5.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.04"
$ws.Range("E34").Value = "  +5.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  +4.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.83"
$ws.Range("E37").Value = "  +8.61%  "
$ws.Range("E38").Value = "  +12.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +8.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.33"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "44.57"
$ws.Range("E41").Value = "  +12.25%  "
$ws.Range("E42").Value = "  +8.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.303"
$ws.Range("E43").Value = "  +13.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.49"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "397.05"
$ws.Range("E45").Value = "  +15.81%  "
$ws.Range("D46").Value = "2.802.05"
$ws.Range("E46").Value = "  +4.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0355"
$ws.Range("E47").Value = "  +5.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.29"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.73"
$ws.Range("E50").Value = "  +10.18%  "
$ws.Range("E51").Value = "  +3.99%  "
